$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.002.67"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.616.09"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "3.011.62"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "2.613.47"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "46.256.16"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "292.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.54%  "
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  +2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "160.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.91%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0844"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("E39").Value = "  +4.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.124"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  +3.34%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.79%  "
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("D46").Value = "2.116.75"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("E47").Value = "  +6.20%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "109.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("E51").Value = "  +0.50%  "
